$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row above the first period row (row 16) ---
# This shifts all period rows (old 16-62) and the signature block (old 67-68)
# down by one, matching the diff (new period rows 16-63, signature rows 68-69).
$ws.Rows("16:16").Insert(-4121)

# Fix up the formatting of the newly inserted row 16: it should look exactly
# like the other "normal" period rows (e.g. row 17), not like the blank
# default that Insert() produces.
$ws.Range("B17:J17").Copy()
$ws.Range("B16:J16").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Rewrite the period list (newest first, plus the new 2507 period) ---
# The worker/employee/amount columns stay the same for every row; only the
# "Periodo Mora" column is re-sequenced, newest period at the top (row 16)
# down to the oldest period (2108) at the bottom (row 63, the bordered
# "last row" of the table).
$periods = @(
    "2507","2506","2505","2504","2503","2502","2501",
    "2412","2411","2410","2409","2408","2407","2406","2405","2404","2403","2402","2401",
    "2312","2311","2310","2309","2308","2307","2306","2305","2304","2303","2302","2301",
    "2212","2211","2210","2209","2208","2207","2206","2205","2204","2203","2202","2201",
    "2112","2111","2110","2109","2108"
)

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Range("B$row").Value = "CC"
    $ws.Range("C$row").Value = "9098406"
    $ws.Range("D$row").Value = "ALEXANDER JULLIO BRU"
    $ws.Range("E$row").Value = $periods[$i]
    $ws.Range("F$row").Value = 47640
    $ws.Range("G$row").Value = 1191000
}

# --- Update the summary figures ---
# VALOR MORA total (E11) grows by the new period's mora value.
$ws.Range("E11").Value = 2286720

# Cant. Periodos (F13) increments from 47 to 48 with the new period row.
$ws.Range("F13").Value = 48
